$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C41").Value = 0.1
$ws.Range("C42").Value = 0.35
$ws.Range("C43").Value = 1
$ws.Range("C44").Value = 0.28000000000000003
$ws.Range("C45").Value = 0.2

$ws.Range("J33").Select()
